{"js": "// Lattice-multiplication exercise sheet: regenerate the 15 practice\n// problems (5 rows x 3 cols) with new operands / digit placeholders.\n// Each table cell holds ONE paragraph / run made of five lines\n// (joined by manual line breaks -> \"\\u000b\" in Office.js text) :\n//   1) \"A x B\"\n//   2) \"  b0    b1\"   (digits of B spaced out)\n//   3) \"  ----\"\n//   4) \"a0|    |\"     (first digit of A)\n//   5) \"a1|    |\"     (second digit of A)\n// The edit only swaps the text content per cell; run formatting\n// (font size 32) is left untouched by replacing the paragraph text\n// in place rather than re-creating runs.\n\nconst cellUpdates = [\n  { row: 0, col: 0, before: \"71 x 50\\u000b  5    0\\u000b  ----\\u000b7|    |\\u000b1|    |\", after: \"52 x 95\\u000b  9    5\\u000b  ----\\u000b5|    |\\u000b2|    |\" },\n  { row: 0, col: 1, before: \"55 x 23\\u000b  2    3\\u000b  ----\\u000b5|    |\\u000b5|    |\", after: \"72 x 99\\u000b  9    9\\u000b  ----\\u000b7|    |\\u000b2|    |\" },\n  { row: 0, col: 2, before: \"22 x 93\\u000b  9    3\\u000b  ----\\u000b2|    |\\u000b2|    |\", after: \"32 x 16\\u000b  1    6\\u000b  ----\\u000b3|    |\\u000b2|    |\" },\n  { row: 1, col: 0, before: \"14 x 24\\u000b  2    4\\u000b  ----\\u000b1|    |\\u000b4|    |\", after: \"14 x 99\\u000b  9    9\\u000b  ----\\u000b1|    |\\u000b4|    |\" },\n  { row: 1, col: 1, before: \"32 x 61\\u000b  6    1\\u000b  ----\\u000b3|    |\\u000b2|    |\", after: \"93 x 84\\u000b  8    4\\u000b  ----\\u000b9|    |\\u000b3|    |\" },\n  { row: 1, col: 2, before: \"99 x 13\\u000b  1    3\\u000b  ----\\u000b9|    |\\u000b9|    |\", after: \"67 x 90\\u000b  9    0\\u000b  ----\\u000b6|    |\\u000b7|    |\" },\n  { row: 2, col: 0, before: \"87 x 69\\u000b  6    9\\u000b  ----\\u000b8|    |\\u000b7|    |\", after: \"94 x 90\\u000b  9    0\\u000b  ----\\u000b9|    |\\u000b4|    |\" },\n  { row: 2, col: 1, before: \"56 x 80\\u000b  8    0\\u000b  ----\\u000b5|    |\\u000b6|    |\", after: \"67 x 60\\u000b  6    0\\u000b  ----\\u000b6|    |\\u000b7|    |\" },\n  { row: 2, col: 2, before: \"64 x 40\\u000b  4    0\\u000b  ----\\u000b6|    |\\u000b4|    |\", after: \"83 x 72\\u000b  7    2\\u000b  ----\\u000b8|    |\\u000b3|    |\" },\n  { row: 3, col: 0, before: \"48 x 39\\u000b  3    9\\u000b  ----\\u000b4|    |\\u000b8|    |\", after: \"33 x 68\\u000b  6    8\\u000b  ----\\u000b3|    |\\u000b3|    |\" },\n  { row: 3, col: 1, before: \"42 x 58\\u000b  5    8\\u000b  ----\\u000b4|    |\\u000b2|    |\", after: \"69 x 15\\u000b  1    5\\u000b  ----\\u000b6|    |\\u000b9|    |\" },\n  { row: 3, col: 2, before: \"80 x 82\\u000b  8    2\\u000b  ----\\u000b8|    |\\u000b0|    |\", after: \"26 x 59\\u000b  5    9\\u000b  ----\\u000b2|    |\\u000b6|    |\" },\n  { row: 4, col: 0, before: \"10 x 91\\u000b  9    1\\u000b  ----\\u000b1|    |\\u000b0|    |\", after: \"45 x 51\\u000b  5    1\\u000b  ----\\u000b4|    |\\u000b5|    |\" },\n  { row: 4, col: 1, before: \"65 x 33\\u000b  3    3\\u000b  ----\\u000b6|    |\\u000b5|    |\", after: \"89 x 19\\u000b  1    9\\u000b  ----\\u000b8|    |\\u000b9|    |\" },\n  { row: 4, col: 2, before: \"40 x 47\\u000b  4    7\\u000b  ----\\u000b4|    |\\u000b0|    |\", after: \"72 x 57\\u000b  5    7\\u000b  ----\\u000b7|    |\\u000b2|    |\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document\");\n}\nconst table = tables.items[0];\n\n// First pass: load every target paragraph's current text.\nconst paras = cellUpdates.map((u) => {\n  const cell = table.getCell(u.row, u.col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  return para;\n});\nawait context.sync();\n\n// Second pass: replace the whole paragraph's text in one shot so the\n// single existing run (and its rPr, e.g. sz=32) is reused instead of\n// being replaced by new default-formatted runs. Only touch cells whose\n// current text still matches the expected \"before\" snapshot.\nfor (let i = 0; i < cellUpdates.length; i++) {\n  const u = cellUpdates[i];\n  const para = paras[i];\n  if (para.text === u.before) {\n    para.insertText(u.after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication exercise sheet: regenerate the 15 practice\n# problems (5 rows x 3 cols) with new operands / digit placeholders.\n# Each table cell holds ONE paragraph / run made of five lines\n# (joined by manual line breaks, char 11 = vertical-tab, which is how\n# Word represents <w:br/> inside Range.Text) :\n#   1) \"A x B\"\n#   2) \"  b0    b1\"   (digits of B spaced out)\n#   3) \"  ----\"\n#   4) \"a0|    |\"     (first digit of A)\n#   5) \"a1|    |\"     (second digit of A)\n# The edit only swaps the text content per cell; assigning Range.Text\n# on the existing cell range reuses the current run (and its rPr,\n# e.g. sz=32) instead of inserting new default-formatted runs.\n\n$cellUpdates = @(\n    [PSCustomObject]@{ Row = 1; Col = 1; Before = @(\"71 x 50\", \"  5    0\", \"  ----\", \"7|    |\", \"1|    |\"); After = @(\"52 x 95\", \"  9    5\", \"  ----\", \"5|    |\", \"2|    |\") }\n    [PSCustomObject]@{ Row = 1; Col = 2; Before = @(\"55 x 23\", \"  2    3\", \"  ----\", \"5|    |\", \"5|    |\"); After = @(\"72 x 99\", \"  9    9\", \"  ----\", \"7|    |\", \"2|    |\") }\n    [PSCustomObject]@{ Row = 1; Col = 3; Before = @(\"22 x 93\", \"  9    3\", \"  ----\", \"2|    |\", \"2|    |\"); After = @(\"32 x 16\", \"  1    6\", \"  ----\", \"3|    |\", \"2|    |\") }\n    [PSCustomObject]@{ Row = 2; Col = 1; Before = @(\"14 x 24\", \"  2    4\", \"  ----\", \"1|    |\", \"4|    |\"); After = @(\"14 x 99\", \"  9    9\", \"  ----\", \"1|    |\", \"4|    |\") }\n    [PSCustomObject]@{ Row = 2; Col = 2; Before = @(\"32 x 61\", \"  6    1\", \"  ----\", \"3|    |\", \"2|    |\"); After = @(\"93 x 84\", \"  8    4\", \"  ----\", \"9|    |\", \"3|    |\") }\n    [PSCustomObject]@{ Row = 2; Col = 3; Before = @(\"99 x 13\", \"  1    3\", \"  ----\", \"9|    |\", \"9|    |\"); After = @(\"67 x 90\", \"  9    0\", \"  ----\", \"6|    |\", \"7|    |\") }\n    [PSCustomObject]@{ Row = 3; Col = 1; Before = @(\"87 x 69\", \"  6    9\", \"  ----\", \"8|    |\", \"7|    |\"); After = @(\"94 x 90\", \"  9    0\", \"  ----\", \"9|    |\", \"4|    |\") }\n    [PSCustomObject]@{ Row = 3; Col = 2; Before = @(\"56 x 80\", \"  8    0\", \"  ----\", \"5|    |\", \"6|    |\"); After = @(\"67 x 60\", \"  6    0\", \"  ----\", \"6|    |\", \"7|    |\") }\n    [PSCustomObject]@{ Row = 3; Col = 3; Before = @(\"64 x 40\", \"  4    0\", \"  ----\", \"6|    |\", \"4|    |\"); After = @(\"83 x 72\", \"  7    2\", \"  ----\", \"8|    |\", \"3|    |\") }\n    [PSCustomObject]@{ Row = 4; Col = 1; Before = @(\"48 x 39\", \"  3    9\", \"  ----\", \"4|    |\", \"8|    |\"); After = @(\"33 x 68\", \"  6    8\", \"  ----\", \"3|    |\", \"3|    |\") }\n    [PSCustomObject]@{ Row = 4; Col = 2; Before = @(\"42 x 58\", \"  5    8\", \"  ----\", \"4|    |\", \"2|    |\"); After = @(\"69 x 15\", \"  1    5\", \"  ----\", \"6|    |\", \"9|    |\") }\n    [PSCustomObject]@{ Row = 4; Col = 3; Before = @(\"80 x 82\", \"  8    2\", \"  ----\", \"8|    |\", \"0|    |\"); After = @(\"26 x 59\", \"  5    9\", \"  ----\", \"2|    |\", \"6|    |\") }\n    [PSCustomObject]@{ Row = 5; Col = 1; Before = @(\"10 x 91\", \"  9    1\", \"  ----\", \"1|    |\", \"0|    |\"); After = @(\"45 x 51\", \"  5    1\", \"  ----\", \"4|    |\", \"5|    |\") }\n    [PSCustomObject]@{ Row = 5; Col = 2; Before = @(\"65 x 33\", \"  3    3\", \"  ----\", \"6|    |\", \"5|    |\"); After = @(\"89 x 19\", \"  1    9\", \"  ----\", \"8|    |\", \"9|    |\") }\n    [PSCustomObject]@{ Row = 5; Col = 3; Before = @(\"40 x 47\", \"  4    7\", \"  ----\", \"4|    |\", \"0|    |\"); After = @(\"72 x 57\", \"  5    7\", \"  ----\", \"7|    |\", \"2|    |\") }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($u in $cellUpdates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    $expectedBefore = [string]::Join([char]11, $u.Before)\n    # Range.Text for a cell includes the trailing end-of-paragraph (CR)\n    # and end-of-cell marker characters; strip them before comparing.\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -eq $expectedBefore) {\n        $newText = [string]::Join([char]11, $u.After)\n        $cell.Range.Text = $newText\n    }\n}\n"}
